$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.511.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "'1.621.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("D5").Value = "'211.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'23.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "

$ws.Range("E9").Value = "  +1.63%  "

$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("D12").Value = "'1.851.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").Value = "'1.625.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").Value = "'0.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.39%  "

$ws.Range("D16").Value = "'65.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").Value = "'27.492.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "'229.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").Value = "'7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("E22").Value = "  +3.35%  "

$ws.Range("D23").Value = "'4.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("E24").Value = "  +7.66%  "

$ws.Range("D25").Value = "'149.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").Value = "'6.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("E27").Value = "  -0.99%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").Value = "'1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("D31").Value = "'0.0484"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").Value = "'1.464.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").Value = "'0.942"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.98%  "

$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").Value = "  -3.29%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  -2.40%  "

$ws.Range("E43").Value = "  -5.47%  "

$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("D46").Value = "'5.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.09%  "

$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("D48").Value = "'1.761.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("D49").Value = "'87.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").Value = "'0.0995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
